$wb = $excel.ActiveWorkbook

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 215921.22
$ws.Range("I15").Value = 215921.22
$ws.Range("K15").Value = 647763.66
$ws.Range("M15").Value = -647594.66

# ALC row 39
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 900.125
$ws.Range("J39").Value = 789.6667
$ws.Range("L39").Value = 2369.0001
$ws.Range("N39").Value = -2961.0001

# ALC row 42
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 736.6667
$ws.Range("J42").Value = 178.1
$ws.Range("L42").Value = 534.3
$ws.Range("N42").Value = -994.3

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1200.5
$ws.Range("I43").Value = 1200.5
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1200.5
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -1131.5
$ws.Range("N43").ClearContents()

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 330573.75
$ws.Range("I98").Value = 533731.4
$ws.Range("J98").Value = 2396
$ws.Range("K98").Value = 533731.4
$ws.Range("L98").Value = 2396
$ws.Range("M98").Value = -532233.4
$ws.Range("N98").Value = -5392

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 330573.75
$ws.Range("I122").Value = 533731.4
$ws.Range("J122").Value = 2396
$ws.Range("K122").Value = 1601194.2
$ws.Range("L122").Value = 7188
$ws.Range("M122").Value = -1598744.2
$ws.Range("N122").Value = -12088

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1034.4375
$ws.Range("I135").Value = 933.6829
$ws.Range("J135").Value = 1624.5714
$ws.Range("K135").Value = 8403.1461
$ws.Range("L135").Value = 14621.1426
$ws.Range("M135").Value = -5868.1461
$ws.Range("N135").Value = -19691.1426

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6700161
$ws.Range("I138").Value = 4636455
$ws.Range("J138").Value = 7095338.5
$ws.Range("K138").Value = 13909365
$ws.Range("L138").Value = 21286015.5
$ws.Range("M138").Value = -13904225
$ws.Range("N138").Value = -21296295.5

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5370.212
$ws.Range("I74").Value = 1618.3636
$ws.Range("J74").Value = 12873.909
$ws.Range("K74").Value = 1618.3636
$ws.Range("L74").Value = 12873.909
$ws.Range("M74").Value = -744.3635999999999
$ws.Range("N74").Value = -14621.909

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5370.212
$ws.Range("I77").Value = 1618.3636
$ws.Range("J77").Value = 12873.909
$ws.Range("K77").Value = 8091.817999999999
$ws.Range("L77").Value = 64369.545
$ws.Range("M77").Value = -3723.817999999999
$ws.Range("N77").Value = -73105.545

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3065.2104
$ws.Range("I132").Value = 2385.6428
$ws.Range("J132").Value = 4968
$ws.Range("K132").Value = 7156.928400000001
$ws.Range("L132").Value = 14904
$ws.Range("M132").Value = -4626.928400000001
$ws.Range("N132").Value = -19964

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 25003332
$ws.Range("I134").Value = 43480332
$ws.Range("J134").Value = 5041.2354
$ws.Range("K134").Value = 130440996
$ws.Range("L134").Value = 15123.7062
$ws.Range("M134").Value = -130438461
$ws.Range("N134").Value = -20193.7062

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2450
$ws.Range("I16").Value = 1720
$ws.Range("J16").Value = 3666.6667
$ws.Range("K16").Value = 1720
$ws.Range("L16").Value = 3666.6667
$ws.Range("M16").Value = -1433
$ws.Range("N16").Value = -4240.6667

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 284.875
$ws.Range("I107").Value = 202.28572
$ws.Range("J107").Value = 318.88235
$ws.Range("K107").Value = 202.28572
$ws.Range("L107").Value = 318.88235
$ws.Range("M107").Value = 1717.71428
$ws.Range("N107").Value = -4158.88235

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2450
$ws.Range("I113").Value = 1720
$ws.Range("J113").Value = 3666.6667
$ws.Range("K113").Value = 1720
$ws.Range("L113").Value = 3666.6667
$ws.Range("M113").Value = 450
$ws.Range("N113").Value = -8006.6667

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1604.4138
$ws.Range("I122").Value = 979.2353000000001
$ws.Range("J122").Value = 2490.0833
$ws.Range("K122").Value = 2937.7059
$ws.Range("L122").Value = 7470.249899999999
$ws.Range("M122").Value = -487.7058999999999
$ws.Range("N122").Value = -12370.2499

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4417.2
$ws.Range("I132").Value = 4510.6665
$ws.Range("J132").Value = 4377.143
$ws.Range("K132").Value = 13531.9995
$ws.Range("L132").Value = 13131.429
$ws.Range("M132").Value = -11001.9995
$ws.Range("N132").Value = -18191.429

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2248.3262
$ws.Range("I134").Value = 1224.8918
$ws.Range("K134").Value = 3674.6754
$ws.Range("M134").Value = -1139.6754

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1171.091
$ws.Range("I5").Value = 468.5
$ws.Range("K5").Value = 1405.5
$ws.Range("M5").Value = -1293.5

# CUL row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1100
$ws.Range("J80").Value = 1150
$ws.Range("L80").Value = 3450
$ws.Range("N80").Value = -5322

# CUL row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 1100
$ws.Range("J83").Value = 1150
$ws.Range("L83").Value = 10350
$ws.Range("N83").Value = -19710

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 583.8570999999999
$ws.Range("I122").Value = 267.53845
$ws.Range("J122").Value = 1097.875
$ws.Range("K122").Value = 2407.84605
$ws.Range("L122").Value = 9880.875
$ws.Range("M122").Value = 42.1539499999999
$ws.Range("N122").Value = -14780.875

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1171.091
$ws.Range("I135").Value = 468.5
$ws.Range("K135").Value = 4216.5
$ws.Range("M135").Value = -1681.5

# GSM row 42
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 14263
$ws.Range("I42").Value = 14263
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 14263
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -13778
$ws.Range("N42").ClearContents()

# GSM row 115
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H115").Value = 14263
$ws.Range("I115").Value = 14263
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 14263
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -13088
$ws.Range("N115").ClearContents()

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1112401.6
$ws.Range("I122").Value = 1588429.8
$ws.Range("J122").Value = 1669.3334
$ws.Range("K122").Value = 4765289.4
$ws.Range("L122").Value = 5008.0002
$ws.Range("M122").Value = -4762839.4
$ws.Range("N122").Value = -9908.0002

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3580.946
$ws.Range("I132").Value = 3532.76
$ws.Range("J132").Value = 3681.3333
$ws.Range("K132").Value = 10598.28
$ws.Range("L132").Value = 11043.9999
$ws.Range("M132").Value = -8068.280000000001
$ws.Range("N132").Value = -16103.9999

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3089.3333
$ws.Range("I7").Value = 1902
$ws.Range("J7").Value = 3428.5715
$ws.Range("K7").Value = 1902
$ws.Range("L7").Value = 3428.5715
$ws.Range("M7").Value = -1790
$ws.Range("N7").Value = -3652.5715

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8029.2666
$ws.Range("I22").Value = 800.6
$ws.Range("J22").Value = 11643.6
$ws.Range("K22").Value = 800.6
$ws.Range("L22").Value = 11643.6
$ws.Range("M22").Value = -505.6
$ws.Range("N22").Value = -12233.6

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 8029.2666
$ws.Range("I27").Value = 800.6
$ws.Range("J27").Value = 11643.6
$ws.Range("K27").Value = 800.6
$ws.Range("L27").Value = 11643.6
$ws.Range("M27").Value = -693.6
$ws.Range("N27").Value = -11857.6

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6369.871
$ws.Range("I61").Value = 6533.522
$ws.Range("J61").Value = 5899.375
$ws.Range("K61").Value = 6533.522
$ws.Range("L61").Value = 5899.375
$ws.Range("M61").Value = -6331.522
$ws.Range("N61").Value = -6303.375

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 6369.871
$ws.Range("I113").Value = 6533.522
$ws.Range("J113").Value = 5899.375
$ws.Range("K113").Value = 6533.522
$ws.Range("L113").Value = 5899.375
$ws.Range("M113").Value = -4363.522
$ws.Range("N113").Value = -10239.375

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3089.3333
$ws.Range("I126").Value = 1902
$ws.Range("J126").Value = 3428.5715
$ws.Range("K126").Value = 5706
$ws.Range("L126").Value = 10285.7145
$ws.Range("M126").Value = -3236
$ws.Range("N126").Value = -15225.7145

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5979.273
$ws.Range("I132").Value = 6009.6665
$ws.Range("K132").Value = 18028.9995
$ws.Range("M132").Value = -15498.9995

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4839.757
$ws.Range("I136").Value = 2890.7144
$ws.Range("J136").Value = 10903.444
$ws.Range("K136").Value = 8672.143199999999
$ws.Range("L136").Value = 32710.332
$ws.Range("M136").Value = -6122.143199999999
$ws.Range("N136").Value = -37810.33199999999

# WVR row 76
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 4500
$ws.Range("J76").Value = 4500
$ws.Range("L76").Value = 4500
$ws.Range("N76").Value = -5130

# WVR row 79
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H79").Value = 4500
$ws.Range("J79").Value = 4500
$ws.Range("L79").Value = 4500
$ws.Range("N79").Value = -6684

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 200680
$ws.Range("I126").Value = 200680
$ws.Range("K126").Value = 602040
$ws.Range("M126").Value = -599570
